# Daily attendance processing - 2025-12-31 20:34:28
# Swap the order of "System" and the email address in the "Recorded By"
# column (G) wherever the cell value is exactly "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7 ("Recorded By")
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
